$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1539.381  # H17
$ws.Cells.Item(17, 10).Value = 1570.3684  # J17
$ws.Cells.Item(17, 12).Value = 4711.1052  # L17
$ws.Cells.Item(17, 14).Value = -5047.1052  # N17
$ws.Cells.Item(33, 8).Value = 280.16666  # H33
$ws.Cells.Item(33, 9).Value = 219.75  # I33
$ws.Cells.Item(33, 10).Value = 401  # J33
$ws.Cells.Item(33, 11).Value = 219.75  # K33
$ws.Cells.Item(33, 12).Value = 401  # L33
$ws.Cells.Item(33, 13).Value = 9.25  # M33
$ws.Cells.Item(33, 14).Value = -859  # N33
$ws.Cells.Item(43, 8).Value = 1762.3334  # H43
$ws.Cells.Item(43, 9).Value = 1727.1  # I43
$ws.Cells.Item(43, 10).Value = 1806.375  # J43
$ws.Cells.Item(43, 11).Value = 1727.1  # K43
$ws.Cells.Item(43, 12).Value = 1806.375  # L43
$ws.Cells.Item(43, 13).Value = -1658.1  # M43
$ws.Cells.Item(43, 14).Value = -1944.375  # N43
$ws.Cells.Item(86, 8).Value = 5740.15  # H86
$ws.Cells.Item(86, 9).Value = 7700.4443  # I86
$ws.Cells.Item(86, 11).Value = 7700.4443  # K86
$ws.Cells.Item(86, 13).Value = -6577.4443  # M86
$ws.Cells.Item(87, 8).Value = 106999.8  # H87
$ws.Cells.Item(87, 10).Value = 121249.75  # J87
$ws.Cells.Item(87, 12).Value = 121249.75  # L87
$ws.Cells.Item(87, 14).Value = -123745.75  # N87
$ws.Cells.Item(89, 8).Value = 5740.15  # H89
$ws.Cells.Item(89, 9).Value = 7700.4443  # I89
$ws.Cells.Item(89, 11).Value = 38502.2215  # K89
$ws.Cells.Item(89, 13).Value = -32886.2215  # M89
$ws.Cells.Item(90, 8).Value = 106999.8  # H90
$ws.Cells.Item(90, 10).Value = 121249.75  # J90
$ws.Cells.Item(90, 12).Value = 363749.25  # L90
$ws.Cells.Item(90, 14).Value = -376229.25  # N90
$ws.Cells.Item(99, 8).Value = 236.63158  # H99
$ws.Cells.Item(99, 9).Value = 194.22223  # I99
$ws.Cells.Item(99, 11).Value = 582.66669  # K99
$ws.Cells.Item(99, 13).Value = 915.33331  # M99
$ws.Cells.Item(101, 8).Value = 310.4  # H101
$ws.Cells.Item(101, 9).Value = 326.75  # I101
$ws.Cells.Item(101, 10).Value = 245  # J101
$ws.Cells.Item(101, 11).Value = 980.25  # K101
$ws.Cells.Item(101, 12).Value = 735  # L101
$ws.Cells.Item(101, 13).Value = 641.75  # M101
$ws.Cells.Item(101, 14).Value = -3979  # N101
$ws.Cells.Item(113, 8).Value = 3646.4546  # H113
$ws.Cells.Item(113, 9).Value = 2599  # I113
$ws.Cells.Item(113, 11).Value = 2599  # K113
$ws.Cells.Item(113, 13).Value = 655  # M113
$ws.Cells.Item(118, 8).Value = 701.4  # H118
$ws.Cells.Item(118, 9).Value = 701.875  # I118
$ws.Cells.Item(118, 11).Value = 2105.625  # K118
$ws.Cells.Item(118, 13).Value = -448.625  # M118
$ws.Cells.Item(121, 8).Value = 2699  # H121
$ws.Cells.Item(121, 10).Value = 2699  # J121
$ws.Cells.Item(121, 12).Value = 8097  # L121
$ws.Cells.Item(121, 14).Value = -11591  # N121
$ws.Cells.Item(138, 8).Value = 2859.9375  # H138
$ws.Cells.Item(138, 9).Value = 2368.0386  # I138
$ws.Cells.Item(138, 10).Value = 3196.5  # J138
$ws.Cells.Item(138, 11).Value = 7104.1158  # K138
$ws.Cells.Item(138, 12).Value = 9589.5  # L138
$ws.Cells.Item(138, 13).Value = -1964.1158  # M138
$ws.Cells.Item(138, 14).Value = -19869.5  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2200.9062  # H61
$ws.Cells.Item(61, 9).Value = 1759.6316  # I61
$ws.Cells.Item(61, 10).Value = 2845.8462  # J61
$ws.Cells.Item(61, 11).Value = 1759.6316  # K61
$ws.Cells.Item(61, 12).Value = 2845.8462  # L61
$ws.Cells.Item(61, 13).Value = -1547.6316  # M61
$ws.Cells.Item(61, 14).Value = -3269.8462  # N61
$ws.Cells.Item(122, 8).Value = 2684.9583  # H122
$ws.Cells.Item(122, 10).Value = 1999.6666  # J122
$ws.Cells.Item(122, 12).Value = 5998.9998  # L122
$ws.Cells.Item(122, 14).Value = -10898.9998  # N122
$ws.Cells.Item(136, 8).Value = 2200.9062  # H136
$ws.Cells.Item(136, 9).Value = 1759.6316  # I136
$ws.Cells.Item(136, 10).Value = 2845.8462  # J136
$ws.Cells.Item(136, 11).Value = 5278.8948  # K136
$ws.Cells.Item(136, 12).Value = 8537.5386  # L136
$ws.Cells.Item(136, 13).Value = -2728.8948  # M136
$ws.Cells.Item(136, 14).Value = -13637.5386  # N136

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 361.63635  # H7
$ws.Cells.Item(7, 9).Value = 177.5  # I7
$ws.Cells.Item(7, 11).Value = 177.5  # K7
$ws.Cells.Item(7, 13).Value = -64.5  # M7
$ws.Cells.Item(16, 8).Value = 1253.0834  # H16
$ws.Cells.Item(16, 9).Value = 395.33334  # I16
$ws.Cells.Item(16, 10).Value = 2110.8333  # J16
$ws.Cells.Item(16, 11).Value = 395.33334  # K16
$ws.Cells.Item(16, 12).Value = 2110.8333  # L16
$ws.Cells.Item(16, 13).Value = -108.33334  # M16
$ws.Cells.Item(16, 14).Value = -2684.8333  # N16
$ws.Cells.Item(31, 8).Value = 2554378.8  # H31
$ws.Cells.Item(31, 9).Value = 2474.3333  # I31
$ws.Cells.Item(31, 11).Value = 2474.3333  # K31
$ws.Cells.Item(31, 13).Value = -2179.3333  # M31
$ws.Cells.Item(34, 8).Value = 2554378.8  # H34
$ws.Cells.Item(34, 9).Value = 2474.3333  # I34
$ws.Cells.Item(34, 11).Value = 2474.3333  # K34
$ws.Cells.Item(34, 13).Value = -2272.3333  # M34
$ws.Cells.Item(56, 8).Value = 0  # H56
$ws.Cells.Item(56, 9).Value = 0  # I56
$ws.Cells.Item(56, 11).Value = 0  # K56
$ws.Cells.Item(56, 13).ClearContents()  # M56
$ws.Cells.Item(58, 8).Value = 1801.9642  # H58
$ws.Cells.Item(58, 9).Value = 1233.1  # I58
$ws.Cells.Item(58, 11).Value = 1233.1  # K58
$ws.Cells.Item(58, 13).Value = -1030.1  # M58
$ws.Cells.Item(99, 8).Value = 3998.25  # H99
$ws.Cells.Item(99, 9).Value = 1996.5  # I99
$ws.Cells.Item(99, 11).Value = 1996.5  # K99
$ws.Cells.Item(99, 13).Value = -498.5  # M99
$ws.Cells.Item(105, 8).Value = 1421.2632  # H105
$ws.Cells.Item(105, 9).Value = 1035.6  # I105
$ws.Cells.Item(105, 11).Value = 1035.6  # K105
$ws.Cells.Item(105, 13).Value = 711.4000000000001  # M105
$ws.Cells.Item(113, 8).Value = 1253.0834  # H113
$ws.Cells.Item(113, 9).Value = 395.33334  # I113
$ws.Cells.Item(113, 10).Value = 2110.8333  # J113
$ws.Cells.Item(113, 11).Value = 395.33334  # K113
$ws.Cells.Item(113, 12).Value = 2110.8333  # L113
$ws.Cells.Item(113, 13).Value = 1774.66666  # M113
$ws.Cells.Item(113, 14).Value = -6450.8333  # N113
$ws.Cells.Item(126, 8).Value = 3998.25  # H126
$ws.Cells.Item(126, 9).Value = 1996.5  # I126
$ws.Cells.Item(126, 11).Value = 5989.5  # K126
$ws.Cells.Item(126, 13).Value = -3519.5  # M126
$ws.Cells.Item(132, 8).Value = 20840402  # H132
$ws.Cells.Item(132, 9).Value = 4899.3335  # I132
$ws.Cells.Item(132, 11).Value = 14698.0005  # K132
$ws.Cells.Item(132, 13).Value = -12168.0005  # M132
$ws.Cells.Item(134, 8).Value = 2262.6743  # H134
$ws.Cells.Item(134, 9).Value = 2008.1621  # I134
$ws.Cells.Item(134, 11).Value = 6024.4863  # K134
$ws.Cells.Item(134, 13).Value = -3489.4863  # M134
$ws.Cells.Item(136, 8).Value = 1801.9642  # H136
$ws.Cells.Item(136, 9).Value = 1233.1  # I136
$ws.Cells.Item(136, 11).Value = 3699.3  # K136
$ws.Cells.Item(136, 13).Value = -1149.3  # M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 7749.1816  # H23
$ws.Cells.Item(23, 9).Value = 291  # I23
$ws.Cells.Item(23, 11).Value = 873  # K23
$ws.Cells.Item(23, 13).Value = -638  # M23
$ws.Cells.Item(93, 8).Value = 10000  # H93
$ws.Cells.Item(93, 10).Value = 10000  # J93
$ws.Cells.Item(93, 12).Value = 30000  # L93
$ws.Cells.Item(93, 14).Value = -33744  # N93
$ws.Cells.Item(141, 8).Value = 23391.5  # H141
$ws.Cells.Item(141, 9).Value = 6783  # I141
$ws.Cells.Item(141, 11).Value = 20349  # K141
$ws.Cells.Item(141, 13).Value = -15169  # M141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 50003444  # H80
$ws.Cells.Item(80, 9).Value = 83336456  # I80
$ws.Cells.Item(80, 11).Value = 83336456  # K80
$ws.Cells.Item(80, 13).Value = -83335458  # M80
$ws.Cells.Item(83, 8).Value = 50003444  # H83
$ws.Cells.Item(83, 9).Value = 83336456  # I83
$ws.Cells.Item(83, 11).Value = 416682280  # K83
$ws.Cells.Item(83, 13).Value = -416677288  # M83
$ws.Cells.Item(96, 8).Value = 29193.166  # H96
$ws.Cells.Item(96, 10).Value = 29193.166  # J96
$ws.Cells.Item(96, 12).Value = 29193.166  # L96
$ws.Cells.Item(96, 14).Value = -34685.166  # N96
$ws.Cells.Item(126, 8).Value = 11828.444  # H126
$ws.Cells.Item(126, 9).Value = 2400  # I126
$ws.Cells.Item(126, 11).Value = 7200  # K126
$ws.Cells.Item(126, 13).Value = -4730  # M126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 900.53845  # H55
$ws.Cells.Item(55, 9).Value = 668.1111  # I55
$ws.Cells.Item(55, 10).Value = 1423.5  # J55
$ws.Cells.Item(55, 11).Value = 668.1111  # K55
$ws.Cells.Item(55, 12).Value = 1423.5  # L55
$ws.Cells.Item(55, 13).Value = -495.1111  # M55
$ws.Cells.Item(55, 14).Value = -1769.5  # N55
$ws.Cells.Item(122, 8).Value = 7698.8647  # H122
$ws.Cells.Item(122, 9).Value = 5187.885  # I122
$ws.Cells.Item(122, 10).Value = 13633.909  # J122
$ws.Cells.Item(122, 11).Value = 15563.655  # K122
$ws.Cells.Item(122, 12).Value = 40901.727  # L122
$ws.Cells.Item(122, 13).Value = -13113.655  # M122
$ws.Cells.Item(122, 14).Value = -45801.727  # N122
$ws.Cells.Item(136, 8).Value = 2781.6  # H136
$ws.Cells.Item(136, 9).Value = 2781.6  # I136
$ws.Cells.Item(136, 11).Value = 8344.799999999999  # K136
$ws.Cells.Item(136, 13).Value = -5794.799999999999  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1082.1666  # H107
$ws.Cells.Item(107, 10).Value = 1771  # J107
$ws.Cells.Item(107, 12).Value = 5313  # L107
$ws.Cells.Item(107, 14).Value = -9153  # N107
$ws.Cells.Item(122, 8).Value = 15628553  # H122
$ws.Cells.Item(122, 9).Value = 3756.4167  # I122
$ws.Cells.Item(122, 11).Value = 11269.2501  # K122
$ws.Cells.Item(122, 13).Value = -8819.250100000001  # M122
